$wb = $excel.ActiveWorkbook
$ws6 = $wb.Worksheets.Item("TC6_Employee_Hostel_Allotment")
$ws7 = $wb.Worksheets.Item("TC7_Hostel_Material")

# Insert a new row 3 on TC7_Hostel_Material ("Executed By" / name), shifting
# everything below it down by one row.
$ws7.Rows.Item(3).Insert()

# The insert spreads formatting across the full row width (A:G); trim that
# back down so only A3/B3 carry content, matching the source row layout.
$ws7.Range("C3:F3").Clear()

# Re-apply the original header-style formatting (copied from the row above)
# before writing in the new label/value pair.
$ws7.Range("A2").Copy()
$ws7.Range("A3").PasteSpecial(-4122)
$ws7.Range("B2").Copy()
$ws7.Range("B3").PasteSpecial(-4122)
$ws7.Range("A3").Value = "Executed By"
$ws7.Range("B3").Value = "Aman Kumar Singh"

# Update the view/selection on TC6_Employee_Hostel_Allotment to focus A3:B3
# and drop the old scrolled-down viewport, without leaving that sheet as the
# active tab.
$null = $ws6.Select()
$null = $ws6.Range("A3:B3").Select()

# Restore TC7_Hostel_Material as the active sheet/tab and update its
# selection to E5.
$null = $ws7.Select()
$null = $ws7.Range("E5").Select()
